$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 350-351, shifting existing rows 350-399 down to 352-401.
# This preserves formatting (style s="2" on column D, a date column) via Excel's
# default "format from cell above" insert behavior.
$ws.Range("A350:R351").Insert(-4121)

# Fill in the new row 350 with its data.
$ws.Range("A350").Value2 = 3
$ws.Range("B350").Value2 = "Femacal de La Calera"
$ws.Range("C350").Value2 = "Coquimbo"
$ws.Range("D350").Value2 = 44491
$ws.Range("E350").Value2 = 5
$ws.Range("F350").Value2 = 100112008
$ws.Range("G350").Value2 = "Coliflor"
$ws.Range("H350").Value2 = "Sin especificar"
$ws.Range("I350").Value2 = "Primera"
$ws.Range("J350").Value2 = 2700
$ws.Range("K350").Value2 = 600
$ws.Range("L350").Value2 = 650
$ws.Range("M350").Value2 = 628
$ws.Range("N350").Value2 = "$/unidad"
$ws.Range("O350").Value2 = "Provincia de Quillota"
$ws.Range("P350").Value2 = 628
$ws.Range("Q350").Value2 = 1
$ws.Range("R350").Value2 = "Hortaliza"

# Fill in the new row 351 with its data.
$ws.Range("A351").Value2 = 3
$ws.Range("B351").Value2 = "Femacal de La Calera"
$ws.Range("C351").Value2 = "Coquimbo"
$ws.Range("D351").Value2 = 44491
$ws.Range("E351").Value2 = 5
$ws.Range("F351").Value2 = 100112008
$ws.Range("G351").Value2 = "Coliflor"
$ws.Range("H351").Value2 = "Sin especificar"
$ws.Range("I351").Value2 = "Segunda"
$ws.Range("J351").Value2 = 1600
$ws.Range("K351").Value2 = 550
$ws.Range("L351").Value2 = 550
$ws.Range("M351").Value2 = 550
$ws.Range("N351").Value2 = "$/unidad"
$ws.Range("O351").Value2 = "Provincia de Quillota"
$ws.Range("P351").Value2 = 550
$ws.Range("Q351").Value2 = 1
$ws.Range("R351").Value2 = "Hortaliza"
